$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely (the "com.hamxa.shaynachim / bitcoin guide /
# bittonnir12@gmail.com / nevilgreen@gmail.com" review row). Every row
# below it shifts up by one.
$ws.Rows("2:2").Delete()

# Row deletion does not renumber this engine's worksheet-level Hyperlinks
# collection, so rebuild it to match the shifted rows: wipe the (now
# stale) collection and re-add the mailto: links at their new locations.
$ws.Range("A1").Hyperlinks.Delete()

$targets = @(
    @{Cell="D5";  Mail="shmulmaor2@gmail.com"},
    @{Cell="C6";  Mail="rocketaso@gmail.com"},
    @{Cell="D6";  Mail="armonravid@gmail.com"},
    @{Cell="C8";  Mail="ronoren61@gmail.com"},
    @{Cell="D8";  Mail="nitanoren23@gmail.com"},
    @{Cell="C10"; Mail="danfogel100@gmail.com"},
    @{Cell="D10"; Mail="avishaybar12@gmail.com"},
    @{Cell="C11"; Mail="danfogel100@gmail.com"},
    @{Cell="D11"; Mail="avishaybar12@gmail.com"},
    @{Cell="D12"; Mail="jorjkluni03@gmail.com"}
)

foreach ($t in $targets) {
    $cellRef = $t.Cell
    $mail = $t.Mail
    $col = $cellRef.Substring(0, 1)
    # Untouched same-column cell (row 4) still carries the original,
    # un-hyperlinked "email" cell style - use it to restore formatting
    # after Hyperlinks.Add() mints its own auto-styled font/underline.
    $fmtSource = $ws.Range($col + "4")
    $target = $ws.Range($cellRef)

    $ws.Hyperlinks.Add($target, "mailto:" + $mail, "", "", $mail)

    $fmtSource.Copy()
    $target.PasteSpecial(-4122)
}

# Match the author's post-delete selection (Excel parks the cursor at B2
# after removing the row).
$ws.Range("B2").Select()
